# Sync attendance_reports: normalize "Recorded By" (column G) entries so
# that the literal "System" marker moves from the front of the
# comma-separated list to the end.
#
# This mirrors the exact textual substitutions observed between the
# before/after workbook revisions:
#   "System, backup@backdoor.com, system" -> "backup@backdoor.com, system, System"
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"         -> "backup@backdoor.com, System"
# Values that don't match one of these exact strings (e.g. the
# "System, admin@admin.com" rows, or rows where "System" is already last)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, system, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$col = 7  # column G = "Recorded By"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
